# agrego número de orden de pago en la validación de pagos en SISE
#
# Updates the three "NroSiniestro" claim-number values used for the
# PREPROD rows (B5:B7) of the SISE payment-order validation sheet, and
# leaves the selection on the last touched cell (B8), matching the
# author's final selection state.
#
# NOTE: a leading apostrophe is used so Excel keeps these numeric-looking
# strings as TEXT (preserving the leading zeros / trailing-space padding)
# instead of coercing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "'0420172008486   "
$ws.Range("B6").Value = "'1220170301430"
$ws.Range("B7").Value = "'1120170200937 "

$ws.Range("B8").Select()
